$p = $ppt.ActivePresentation

# The "Q & A" slide (position 35) is replaced with a new "Homework" slide.
# Duplicating it first (then deleting the original) gives the replacement
# slide a fresh internal slide id, matching how PowerPoint mints a new id
# whenever a slide is effectively re-created, while keeping the clean
# placeholder/shape structure of the original.
$old = $p.Slides.Item(35)
$dupRange = $old.Duplicate()
$slide = $dupRange.Item(1)
$old.Delete()

# --- Title placeholder: "Q & A" -> "Homework" ---
$title = $slide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Homework"

# --- Body placeholder: homework link + helpful module references ---
# Paragraph layout (1-based):
#   1: hyperlinked homework repo link
#   2: (blank)
#   3: "The following module references might be helpful:"
#   4: (indented) hyperlinked Enum.html reference
#   5: (indented) hyperlinked Map.html reference
$body = $slide.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "https://github.com/georgiyolovski/elixir-workshop/tree/main/day1/homework`rX`rThe following module references might be helpful:`rhttps://hexdocs.pm/elixir/Enum.html`rhttps://hexdocs.pm/elixir/Map.html"
$body.TextFrame.TextRange.LanguageID = "en-US"

# Turn paragraph 2 into a genuinely empty paragraph (placeholder "X" typed
# then removed so we don't end up with a stray empty run).
$body.TextFrame.TextRange.Paragraphs(2).Characters(1, 1).Delete()

$body.TextFrame.TextRange.Paragraphs(4).IndentLevel = 2
$body.TextFrame.TextRange.Paragraphs(5).IndentLevel = 2

$body.TextFrame.TextRange.Paragraphs(1).ActionSettings.Item(1).Hyperlink.Address = "https://github.com/georgiyolovski/elixir-workshop/tree/main/day1/homework"
$body.TextFrame.TextRange.Paragraphs(4).ActionSettings.Item(1).Hyperlink.Address = "https://hexdocs.pm/elixir/Enum.html"
$body.TextFrame.TextRange.Paragraphs(5).ActionSettings.Item(1).Hyperlink.Address = "https://hexdocs.pm/elixir/Map.html"

# Match the target shape names for the rebuilt slide.
$title.Name = "Title 3"
$body.Name = "Text Placeholder 4"
